$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166006922721863
$ws.Range("B1").Value = 2.435671091079712
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.368738651275635
$ws.Range("E1").Value = 1.235018253326416
